# Scheduled-runner market-price refresh for Lich_Profits sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the derived
# LeveProfit(NQ/HQ) columns (H:N) for the leves whose item prices moved.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 103.85714
$ws.Range("I2").Value = 109
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 109
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = -326

# Row 39
$ws.Range("H39").Value = 1526.125
$ws.Range("I39").Value = 1461.8
$ws.Range("K39").Value = 4385.4
$ws.Range("M39").Value = -4089.4

# Row 62
$ws.Range("H62").Value = 7414.643
$ws.Range("I62").Value = 6991.4546
$ws.Range("K62").Value = 6991.4546
$ws.Range("M62").Value = -6367.4546

# Row 64
$ws.Range("H64").Value = 6060.8335
$ws.Range("I64").Value = 5584.3335
$ws.Range("K64").Value = 5584.3335
$ws.Range("M64").Value = -5336.3335

# Row 65
$ws.Range("H65").Value = 7414.643
$ws.Range("I65").Value = 6991.4546
$ws.Range("K65").Value = 34957.273
$ws.Range("M65").Value = -31837.273

# Row 67
$ws.Range("H67").Value = 6060.8335
$ws.Range("I67").Value = 5584.3335
$ws.Range("K67").Value = 5584.3335
$ws.Range("M67").Value = -4726.3335

# Row 76
$ws.Range("H76").Value = 3702
$ws.Range("I76").Value = 3515.3333
$ws.Range("K76").Value = 3515.3333
$ws.Range("M76").Value = -3200.3333

# Row 79
$ws.Range("H79").Value = 3702
$ws.Range("I79").Value = 3515.3333
$ws.Range("K79").Value = 3515.3333
$ws.Range("M79").Value = -2423.3333

# Row 88
$ws.Range("H88").Value = 2398.5
$ws.Range("I88").Value = 2031.6666
$ws.Range("K88").Value = 2031.6666
$ws.Range("M88").Value = -1625.6666

# Row 91
$ws.Range("H91").Value = 2398.5
$ws.Range("I91").Value = 2031.6666
$ws.Range("K91").Value = 2031.6666
$ws.Range("M91").Value = -627.6666

# Row 137
$ws.Range("H137").Value = 65135.473
$ws.Range("I137").Value = 72680.94
$ws.Range("K137").Value = 218042.82
$ws.Range("M137").Value = -215492.82

# Row 138
$ws.Range("H138").Value = 3230.3
$ws.Range("I138").Value = 2087.12
$ws.Range("J138").Value = 3611.36
$ws.Range("K138").Value = 6261.36
$ws.Range("L138").Value = 10834.08
$ws.Range("M138").Value = -1121.36
$ws.Range("N138").Value = -21114.08


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2357.4443
$ws.Range("I63").Value = 2460.6667
$ws.Range("K63").Value = 2460.6667
$ws.Range("M63").Value = -1774.6667

# Row 66
$ws.Range("H66").Value = 2357.4443
$ws.Range("I66").Value = 2460.6667
$ws.Range("K66").Value = 12303.3335
$ws.Range("M66").Value = -8871.333500000001

# Row 88
$ws.Range("H88").Value = 1600.4736
$ws.Range("I88").Value = 1777.5
$ws.Range("J88").Value = 1518.7693
$ws.Range("K88").Value = 1777.5
$ws.Range("L88").Value = 1518.7693
$ws.Range("M88").Value = -1371.5
$ws.Range("N88").Value = -2330.7693

# Row 91
$ws.Range("H91").Value = 1600.4736
$ws.Range("I91").Value = 1777.5
$ws.Range("J91").Value = 1518.7693
$ws.Range("K91").Value = 1777.5
$ws.Range("L91").Value = 1518.7693
$ws.Range("M91").Value = -373.5
$ws.Range("N91").Value = -4326.7693


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1228.62
$ws.Range("I86").Value = 1145.3077
$ws.Range("J86").Value = 1524
$ws.Range("K86").Value = 1145.3077
$ws.Range("L86").Value = 1524
$ws.Range("M86").Value = -22.30770000000007
$ws.Range("N86").Value = -3770

# Row 89
$ws.Range("H89").Value = 1228.62
$ws.Range("I89").Value = 1145.3077
$ws.Range("J89").Value = 1524
$ws.Range("K89").Value = 5726.538500000001
$ws.Range("L89").Value = 7620
$ws.Range("M89").Value = -110.5385000000006
$ws.Range("N89").Value = -18852

# Row 105
$ws.Range("H105").Value = 3117.6
$ws.Range("I105").Value = 3152.889
$ws.Range("K105").Value = 3152.889
$ws.Range("M105").Value = -1405.889


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2486
$ws.Range("I16").Value = 2380.5715
$ws.Range("J16").Value = 2855
$ws.Range("K16").Value = 2380.5715
$ws.Range("L16").Value = 2855
$ws.Range("M16").Value = -2093.5715
$ws.Range("N16").Value = -3429

# Row 62
$ws.Range("H62").Value = 8321.143
$ws.Range("I62").Value = 9498.9
$ws.Range("J62").Value = 5376.75
$ws.Range("K62").Value = 9498.9
$ws.Range("L62").Value = 5376.75
$ws.Range("M62").Value = -8874.9
$ws.Range("N62").Value = -6624.75

# Row 65
$ws.Range("H65").Value = 8321.143
$ws.Range("I65").Value = 9498.9
$ws.Range("J65").Value = 5376.75
$ws.Range("K65").Value = 47494.5
$ws.Range("L65").Value = 26883.75
$ws.Range("M65").Value = -44374.5
$ws.Range("N65").Value = -33123.75

# Row 113
$ws.Range("H113").Value = 2486
$ws.Range("I113").Value = 2380.5715
$ws.Range("J113").Value = 2855
$ws.Range("K113").Value = 2380.5715
$ws.Range("L113").Value = 2855
$ws.Range("M113").Value = -210.5715
$ws.Range("N113").Value = -7195

# Row 122
$ws.Range("H122").Value = 2312.5334
$ws.Range("I122").Value = 1984.4546
$ws.Range("K122").Value = 5953.3638
$ws.Range("M122").Value = -3503.3638


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 2571.0938
$ws.Range("I12").Value = 457.2
$ws.Range("J12").Value = 3531.9546
$ws.Range("K12").Value = 1371.6
$ws.Range("L12").Value = 10595.8638
$ws.Range("M12").Value = -1198.6
$ws.Range("N12").Value = -10941.8638


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10783.523
$ws.Range("I70").Value = 10351.308
$ws.Range("K70").Value = 10351.308
$ws.Range("M70").Value = -10081.308

# Row 73
$ws.Range("H73").Value = 10783.523
$ws.Range("I73").Value = 10351.308
$ws.Range("K73").Value = 10351.308
$ws.Range("M73").Value = -9415.308000000001

# Row 80
$ws.Range("H80").Value = 3262.04
$ws.Range("J80").Value = 4899.5
$ws.Range("L80").Value = 4899.5
$ws.Range("N80").Value = -6895.5

# Row 83
$ws.Range("H83").Value = 3262.04
$ws.Range("J83").Value = 4899.5
$ws.Range("L83").Value = 24497.5
$ws.Range("N83").Value = -34481.5

# Row 126
$ws.Range("H126").Value = 6777.25
$ws.Range("I126").Value = 5344.5
$ws.Range("K126").Value = 16033.5
$ws.Range("M126").Value = -13563.5


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1886.037
$ws.Range("I16").Value = 542.8182
$ws.Range("J16").Value = 7796.2
$ws.Range("K16").Value = 542.8182
$ws.Range("L16").Value = 7796.2
$ws.Range("M16").Value = -372.8182
$ws.Range("N16").Value = -8136.2


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2957.1765
$ws.Range("I81").Value = 3023.125
$ws.Range("J81").Value = 1902
$ws.Range("K81").Value = 6046.25
$ws.Range("L81").Value = 3804
$ws.Range("M81").Value = -4985.25
$ws.Range("N81").Value = -5926

# Row 84
$ws.Range("H84").Value = 2957.1765
$ws.Range("I84").Value = 3023.125
$ws.Range("J84").Value = 1902
$ws.Range("K84").Value = 30231.25
$ws.Range("L84").Value = 19020
$ws.Range("M84").Value = -24927.25
$ws.Range("N84").Value = -29628
